# Update automatico via Actualizar 02-05-2021 14-21-37
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Correct the D604:D617 timestamps (tiny float re-compute) ---
$ws.Range("D604:D617").Value = 44232.57720748842

# --- 2) Append the new batch of 14 rows (618-631), same 14-service cycle ---
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls  = @(
  "https://www.dataintelligence-group.com/",
  "https://serviciodashboard.azurewebsites.net/",
  "https://powerbi.microsoft.com/es-es/",
  "https://www.dropbox.com/",
  "https://dataintelligence.store/",
  "https://app-data-i.users.earthengine.app/",
  "https://odooutil.azurewebsites.net/",
  "https://filtradordashboard.azurewebsites.net/",
  "https://ide.dataintelligence-group.com/mapstore/#/",
  "https://ide.dataintelligence-group.com/geoserver/web/?0",
  "https://ide.dataintelligence-group.com/",
  "https://rpubs.com/dataintelligence/",
  "https://github.com/Sud-Austral/",
  "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$hlAddress = @(
  "https://www.dataintelligence-group.com/",
  "https://serviciodashboard.azurewebsites.net/",
  "https://powerbi.microsoft.com/es-es/",
  "https://www.dropbox.com/",
  "https://dataintelligence.store/",
  "https://app-data-i.users.earthengine.app/",
  "https://odooutil.azurewebsites.net/",
  "https://filtradordashboard.azurewebsites.net/",
  "https://ide.dataintelligence-group.com/mapstore/",
  "https://ide.dataintelligence-group.com/geoserver/web/?0",
  "https://ide.dataintelligence-group.com/",
  "https://rpubs.com/dataintelligence/",
  "https://github.com/Sud-Austral/",
  "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$hlSub = @("","","","","","","","","/","","","","","")

$startRow = 618
for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = $names[$i]
    $ws.Range("B$r").Value = $urls[$i]
    $ws.Range("C$r").Value = "Disponible"
    $ws.Range("D$r").Value = 44232.59827623806

    # Copy formatting (styles) down from the row above (keeps shared style ids: s=2 hyperlink, s=3 date)
    $ws.Range("A" + ($r - 1) + ":D" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = 0

    # Register the hyperlink relationship for column B
    $ws.Hyperlinks.Add($ws.Range("B$r"), $hlAddress[$i], $hlSub[$i])

    # Re-apply the formatting once more so the Hyperlinks.Add call's own
    # style side-effect doesn't stick on the cell
    $ws.Range("A" + ($r - 1) + ":D" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = 0
}

Write-Host "Added rows 618-631"
